$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "61.028.69"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.663.72"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "530.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.08%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "155.74"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.583"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "6.50"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("E10").Value = "  +4.65%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -0.49%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "3.126.96"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.80%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "61.046.61"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("E16").Value = "  +1.41%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.663.41"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.62%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "354.67"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.69"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.34"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  +0.71%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0₃0858"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.33"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -0.03%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("E31").Value = "  +0.21%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.63"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.59%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "149.96"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +0.22%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.920"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +7.89%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.894"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "37.02"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "305.81"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.16%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.49"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("E43").Value = "  +0.28%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "20.37"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.46%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  +2.39%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "4.91"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "19.34"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.84%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "10.35"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.994.99"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
